$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.426.89'
$ws.Range("E2").Value = '  -1.37%  '

# Row 3
$ws.Range("D3").Value = '2.374.29'
$ws.Range("E3").Value = '  +4.63%  '

# Row 4
$ws.Range("E4").Value = '  +0.12%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.39'
$ws.Range("E5").Value = '  +0.79%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.655'
$ws.Range("E6").Value = '  -0.07%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '71.83'
$ws.Range("E7").Value = '  +12.40%  '

# Row 8
$ws.Range("E8").Value = '  +0.06%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.471'
$ws.Range("E9").Value = '  +3.85%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0978'
$ws.Range("E10").Value = '  -0.68%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '56.86'
$ws.Range("E11").Value = '  -1.76%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '27.20'
$ws.Range("E12").Value = '  +0.99%  '

# Row 13
$ws.Range("D13").Value = '2.727.41'
$ws.Range("E13").Value = '  +4.71%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.106'
$ws.Range("E14").Value = '  -0.04%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.94'
$ws.Range("E15").Value = '  +1.25%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.30'
$ws.Range("E16").Value = '  +1.97%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.855'
$ws.Range("E17").Value = '  +1.45%  '

# Row 18
$ws.Range("D18").Value = '2.371.45'
$ws.Range("E18").Value = '  +4.62%  '

# Row 19
$ws.Range("D19").Value = '43.476.85'
$ws.Range("E19").Value = '  -1.07%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0999'
$ws.Range("E20").Value = '  +0.98%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '74.87'
$ws.Range("E21").Value = '  +0.85%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.34'
$ws.Range("E22").Value = '  +3.44%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '251.50'
$ws.Range("E23").Value = '  +0.09%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.81'
$ws.Range("E24").Value = '  +14.77%  '

# Row 25
$ws.Range("E25").Value = '  -0.06%  '

# Row 26
$ws.Range("E26").Value = '  +1.38%  '

# Row 27
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.27'
$ws.Range("E27").Value = '  +1.56%  '

# Row 28
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.90'
$ws.Range("E28").Value = '  +2.07%  '

# Row 29
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.00'
$ws.Range("E29").Value = '  +0.05%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '174.47'
$ws.Range("E30").Value = '  +0.36%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.52'
$ws.Range("E31").Value = '  +3.93%  '

# Row 32
$ws.Range("E32").Value = '  -7.11%  '

# Row 33
$ws.Range("E33").Value = '  -0.01%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.01'
$ws.Range("E34").Value = '  -0.62%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0691'
$ws.Range("E35").Value = '  +0.70%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.06'
$ws.Range("E36").Value = '  +0.62%  '

# Row 37
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.45'
$ws.Range("E37").Value = '  +6.12%  '

# Row 38
$ws.Range("B38").Value = 'THORChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.59'
$ws.Range("E38").Value = '  +1.69%  '

# Row 39
$ws.Range("E39").Value = '  -1.03%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0257'
$ws.Range("E40").Value = '  +0.22%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.95'
$ws.Range("E41").Value = '  +1.26%  '

# Row 42
$ws.Range("B42").Value = 'BinanceUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.08%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '18.59'
$ws.Range("E43").Value = '  +6.97%  '

# Row 44
$ws.Range("E44").Value = '  +7.29%  '

# Row 45
$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.54'
$ws.Range("E45").Value = '  +3.65%  '

# Row 46
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '100.26'
$ws.Range("E46").Value = '  +1.50%  '

# Row 47
$ws.Range("E47").Value = '  +1.74%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0957'
$ws.Range("E48").Value = '  +0.26%  '

# Row 49
$ws.Range("B49").Value = 'TerraClassic'
$ws.Range("C49").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000207'
$ws.Range("E49").Value = '  -6.71%  '

# Row 50
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '1.444.37'
$ws.Range("E50").Value = '  -0.72%  '

# Row 51
$ws.Range("D51").Value = '2.599.80'
